# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2404   (left-hand "before" columns, A1:J1)
#   *_new -> *_FV2410   (right-hand "after" columns, L1:U1)
# and turn the data range into a proper Excel Table with a frozen header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells in row 1 -----------------------------------
$lastCol = 21  # A..U
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = [string]$cell.Value2
    if ($val.EndsWith("_old")) {
        $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2404"
    } elseif ($val.EndsWith("_new")) {
        $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2410"
    }
}

# --- 2. Turn the used range into an Excel Table -----------------------------
$usedRange = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
